$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 16587.938
$ws.Range("J32").Value = 19027.7
$ws.Range("L32").Value = 19027.7
$ws.Range("N32").Value = -19679.7
$ws.Range("H33").Value = 311.15384
$ws.Range("I33").Value = 171.875
$ws.Range("J33").Value = 534
$ws.Range("K33").Value = 171.875
$ws.Range("L33").Value = 534
$ws.Range("M33").Value = 57.125
$ws.Range("N33").Value = -992
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 20000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 20000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -21944
$ws.Range("H98").Value = 3987.1667
$ws.Range("I98").Value = 1984.875
$ws.Range("K98").Value = 1984.875
$ws.Range("M98").Value = -486.875
$ws.Range("H104").Value = 100
$ws.Range("I104").Value = 100
$ws.Range("K104").Value = 300
$ws.Range("M104").Value = 1447
$ws.Range("H122").Value = 3987.1667
$ws.Range("I122").Value = 1984.875
$ws.Range("K122").Value = 5954.625
$ws.Range("M122").Value = -3504.625
$ws.Range("H127").Value = 2257.6365
$ws.Range("I127").Value = 1397.4
$ws.Range("K127").Value = 4192.200000000001
$ws.Range("M127").Value = 767.7999999999993
$ws.Range("H137").Value = 51783
$ws.Range("I137").Value = 50765.668
$ws.Range("K137").Value = 152297.004
$ws.Range("M137").Value = -149747.004
$ws.Range("H141").Value = 2888.2856
$ws.Range("I141").Value = 2905.75
$ws.Range("K141").Value = 8717.25
$ws.Range("M141").Value = -3537.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H97").Value = 672.1591
$ws.Range("I97").Value = 633.8293
$ws.Range("K97").Value = 633.8293
$ws.Range("M97").Value = -137.8293
$ws.Range("H102").Value = 2568.1765
$ws.Range("I102").Value = 2425.25
$ws.Range("K102").Value = 2425.25
$ws.Range("M102").Value = -803.25
$ws.Range("H122").Value = 1815.8422
$ws.Range("I122").Value = 1592.2858
$ws.Range("K122").Value = 4776.857400000001
$ws.Range("M122").Value = -2326.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2279.8
$ws.Range("I86").Value = 2199.7778
$ws.Range("K86").Value = 2199.7778
$ws.Range("M86").Value = -1076.7778
$ws.Range("H89").Value = 2279.8
$ws.Range("I89").Value = 2199.7778
$ws.Range("K89").Value = 10998.889
$ws.Range("M89").Value = -5382.888999999999
$ws.Range("H99").Value = 1059.25
$ws.Range("I99").Value = 1064.7273
$ws.Range("J99").Value = 999
$ws.Range("K99").Value = 1064.7273
$ws.Range("L99").Value = 999
$ws.Range("M99").Value = 433.2727
$ws.Range("N99").Value = -3995
$ws.Range("H107").Value = 1906.7222
$ws.Range("I107").Value = 1958.3429
$ws.Range("K107").Value = 1958.3429
$ws.Range("M107").Value = -38.3429000000001
$ws.Range("H134").Value = 2601.8
$ws.Range("I134").Value = 1729.1428
$ws.Range("K134").Value = 5187.428400000001
$ws.Range("M134").Value = -2652.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 795
$ws.Range("J22").Value = 833.3333
$ws.Range("L22").Value = 833.3333
$ws.Range("N22").Value = -1533.3333
$ws.Range("H31").Value = 5556559.5
$ws.Range("J31").Value = 1277.25
$ws.Range("L31").Value = 1277.25
$ws.Range("N31").Value = -1867.25
$ws.Range("H34").Value = 5556559.5
$ws.Range("J34").Value = 1277.25
$ws.Range("L34").Value = 1277.25
$ws.Range("N34").Value = -1681.25
$ws.Range("H105").Value = 1545.8235
$ws.Range("I105").Value = 720.7692
$ws.Range("J105").Value = 4227.25
$ws.Range("K105").Value = 720.7692
$ws.Range("L105").Value = 4227.25
$ws.Range("M105").Value = 1026.2308
$ws.Range("N105").Value = -7721.25
$ws.Range("H107").Value = 652.5599999999999
$ws.Range("I107").Value = 442.5263
$ws.Range("J107").Value = 1317.6666
$ws.Range("K107").Value = 442.5263
$ws.Range("L107").Value = 1317.6666
$ws.Range("M107").Value = 1477.4737
$ws.Range("N107").Value = -5157.6666
$ws.Range("H122").Value = 1290.5
$ws.Range("I122").Value = 1338.7333
$ws.Range("K122").Value = 4016.199900000001
$ws.Range("M122").Value = -1566.199900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 104538480
$ws.Range("I4").Value = 52433588
$ws.Range("K4").Value = 157300764
$ws.Range("M4").Value = -157300652
$ws.Range("H9").Value = 750
$ws.Range("I9").Value = 500
$ws.Range("K9").Value = 1500
$ws.Range("M9").Value = -1276
$ws.Range("H88").Value = 5000
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856
$ws.Range("H91").Value = 5000
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964
$ws.Range("H92").Value = 1632.9166
$ws.Range("J92").Value = 2124.5
$ws.Range("L92").Value = 6373.5
$ws.Range("N92").Value = -8869.5
$ws.Range("H97").Value = 861.4545000000001
$ws.Range("J97").Value = 892.3333
$ws.Range("L97").Value = 2676.9999
$ws.Range("N97").Value = -3668.9999
$ws.Range("H121").Value = 81176.35000000001
$ws.Range("I121").Value = 133601.25
$ws.Range("J121").Value = 34576.445
$ws.Range("K121").Value = 400803.75
$ws.Range("L121").Value = 103729.335
$ws.Range("M121").Value = -399493.75
$ws.Range("N121").Value = -106349.335
$ws.Range("H132").Value = 2299.3333
$ws.Range("J132").Value = 2299.3333
$ws.Range("L132").Value = 20693.9997
$ws.Range("N132").Value = -25753.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 1906625
$ws.Range("J44").Value = 1906625
$ws.Range("L44").Value = 1906625
$ws.Range("N44").Value = -1907817
$ws.Range("H70").Value = 5827.1113
$ws.Range("I70").Value = 7283
$ws.Range("J70").Value = 5099.1665
$ws.Range("K70").Value = 7283
$ws.Range("L70").Value = 5099.1665
$ws.Range("M70").Value = -7013
$ws.Range("N70").Value = -5639.1665
$ws.Range("H73").Value = 5827.1113
$ws.Range("I73").Value = 7283
$ws.Range("J73").Value = 5099.1665
$ws.Range("K73").Value = 7283
$ws.Range("L73").Value = 5099.1665
$ws.Range("M73").Value = -6347
$ws.Range("N73").Value = -6971.1665
$ws.Range("H122").Value = 2730.258
$ws.Range("I122").Value = 2549.9167
$ws.Range("J122").Value = 3348.5715
$ws.Range("K122").Value = 7649.750100000001
$ws.Range("L122").Value = 10045.7145
$ws.Range("M122").Value = -5199.750100000001
$ws.Range("N122").Value = -14945.7145
$ws.Range("H126").Value = 3521.9333
$ws.Range("I126").Value = 1925
$ws.Range("J126").Value = 5347
$ws.Range("K126").Value = 5775
$ws.Range("L126").Value = 16041
$ws.Range("M126").Value = -3305
$ws.Range("N126").Value = -20981
$ws.Range("H136").Value = 98965
$ws.Range("J136").Value = 98965
$ws.Range("L136").Value = 296895
$ws.Range("N136").Value = -301995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4120.727
$ws.Range("I7").Value = 4568.25
$ws.Range("J7").Value = 2927.3333
$ws.Range("K7").Value = 4568.25
$ws.Range("L7").Value = 2927.3333
$ws.Range("M7").Value = -4456.25
$ws.Range("N7").Value = -3151.3333
$ws.Range("H16").Value = 470
$ws.Range("I16").Value = 332.66666
$ws.Range("K16").Value = 332.66666
$ws.Range("M16").Value = -162.66666
$ws.Range("H40").Value = 4055.2
$ws.Range("I40").Value = 4071.111
$ws.Range("J40").Value = 3912
$ws.Range("K40").Value = 4071.111
$ws.Range("L40").Value = 3912
$ws.Range("M40").Value = -3935.111
$ws.Range("N40").Value = -4184
$ws.Range("H46").Value = 4268.143
$ws.Range("I46").Value = 1598.2
$ws.Range("J46").Value = 5751.4443
$ws.Range("K46").Value = 1598.2
$ws.Range("L46").Value = 5751.4443
$ws.Range("M46").Value = -1410.2
$ws.Range("N46").Value = -6127.4443
$ws.Range("H100").Value = 2489.0417
$ws.Range("I100").Value = 2401.6843
$ws.Range("K100").Value = 2401.6843
$ws.Range("M100").Value = -1860.6843
$ws.Range("H126").Value = 4120.727
$ws.Range("I126").Value = 4568.25
$ws.Range("J126").Value = 2927.3333
$ws.Range("K126").Value = 13704.75
$ws.Range("L126").Value = 8781.999899999999
$ws.Range("M126").Value = -11234.75
$ws.Range("N126").Value = -13721.9999
$ws.Range("H132").Value = 7388.1113
$ws.Range("I132").Value = 7784.7144
$ws.Range("K132").Value = 23354.1432
$ws.Range("M132").Value = -20824.1432
$ws.Range("H134").Value = 97500
$ws.Range("I134").Value = 75000
$ws.Range("J134").Value = 108750
$ws.Range("K134").Value = 75000
$ws.Range("L134").Value = 108750
$ws.Range("M134").Value = -69930
$ws.Range("N134").Value = -118890
$ws.Range("H136").Value = 5334
$ws.Range("I136").Value = 5334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -13452
$ws.Range("N136").ClearContents()
$ws.Range("H140").Value = 107999.5
$ws.Range("J140").Value = 107999.5
$ws.Range("L140").Value = 107999.5
$ws.Range("N140").Value = -118359.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H141").Value = 116998.11
$ws.Range("I141").Value = 84999
$ws.Range("J141").Value = 120998
$ws.Range("K141").Value = 84999
$ws.Range("L141").Value = 120998
$ws.Range("M141").Value = -79819
$ws.Range("N141").Value = -131358
